$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / changed date) from 45221 (2023-10-22) to
# 45224 (2023-10-25) for all data rows (rows 2 through 54).
for ($r = 2; $r -le 54; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
